$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E/F column type-pair text + outcome corrections ---
# (fill colours below mirror the existing Win/Neutral/Loss conditional
#  formatting already used throughout column F: green/orange/red)
$ws.Range("E106").Value = '{''int'', ''any''}'
$ws.Range("E107").Value = 'int'
$ws.Range("F107").Interior.Color = 0x008000
$ws.Range("F107").Value = 'Win'
$ws.Range("E112").Value = '{''Atom'', ''empty''}'
$ws.Range("E113").Value = 'Atom'
$ws.Range("F113").Interior.Color = 0x00A5FF
$ws.Range("F113").Value = 'Neutral'
$ws.Range("E114").Value = '{''Udta'', ''empty''}'
$ws.Range("E115").Value = 'Udta'
$ws.Range("E116").Value = '{''Tuple[any]'', ''empty''}'
$ws.Range("E117").Value = 'Tuple[any]'
$ws.Range("F117").Interior.Color = 0x00A5FF
$ws.Range("F117").Value = 'Neutral'
$ws.Range("E122").Value = '{''VariableAtom'', ''empty''}'
$ws.Range("E123").Value = 'VariableAtom'
$ws.Range("F123").Interior.Color = 0x00A5FF
$ws.Range("F123").Value = 'Neutral'
$ws.Range("E124").Value = '{''VariableAtom'', ''empty''}'
$ws.Range("E125").Value = 'VariableAtom'
$ws.Range("F125").Interior.Color = 0x00A5FF
$ws.Range("F125").Value = 'Neutral'
$ws.Range("E126").Value = '{''VariableAtom'', ''empty''}'
$ws.Range("E127").Value = 'VariableAtom'
$ws.Range("F127").Interior.Color = 0x00A5FF
$ws.Range("F127").Value = 'Neutral'
$ws.Range("E130").Value = '{''Atom'', ''empty''}'
$ws.Range("E131").Value = 'Atom'
$ws.Range("F131").Interior.Color = 0x00A5FF
$ws.Range("F131").Value = 'Neutral'
$ws.Range("E134").Value = '{''Atom'', ''empty''}'
$ws.Range("E135").Value = 'Atom'
$ws.Range("F135").Interior.Color = 0x00A5FF
$ws.Range("F135").Value = 'Neutral'
$ws.Range("E150").Value = '{''Atom'', ''empty''}'
$ws.Range("E151").Value = 'Atom'
$ws.Range("F151").Interior.Color = 0x00A5FF
$ws.Range("F151").Value = 'Neutral'
$ws.Range("E152").Value = '{''Atom'', ''empty''}'
$ws.Range("E153").Value = 'Atom'
$ws.Range("F153").Interior.Color = 0x00A5FF
$ws.Range("F153").Value = 'Neutral'
$ws.Range("E154").Value = '{''Atom'', ''empty''}'
$ws.Range("E155").Value = 'Atom'
$ws.Range("F155").Interior.Color = 0x00A5FF
$ws.Range("F155").Value = 'Neutral'
$ws.Range("E156").Value = '{''CompositeAtom'', ''empty''}'
$ws.Range("E157").Value = 'CompositeAtom'
$ws.Range("F157").Interior.Color = 0x00A5FF
$ws.Range("F157").Value = 'Neutral'
$ws.Range("E477").Value = '{''any'', ''Tuple[None]''}'
$ws.Range("E478").Value = 'any'

# --- Summary block (rows 503-505): fix accuracy computation ---
# Previously PyType Wins/Scalpel Wins were mis-tallied (D503/F503),
# and "Scalpel Accuracy" wrongly divided by PyType Wins (D503) instead of
# the total comparison count (B503), and its label/value sat in C504/D504
# instead of lining up with the other summary rows in E/F.
$ws.Range("D503").Value = 12
$ws.Range("F503").Value = 117

$ws.Range("C504").ClearContents()
$ws.Range("D504").ClearContents()
$ws.Range("E504").Value = "Scalpel Accuracy:"
$ws.Range("F504").Value = 97.6

$ws.Range("E505").Value = "Accuracy vs PyType"
$ws.Range("F505").Value = 975
